# Restored from revision #7b36ba64c53effbf799508ca8248539a2296c88e.TEST
# Author: admin. Type: SAVE.
#
# Update the "Integer max" value of rule R20 on the Rules sheet from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
